$wb = $excel.ActiveWorkbook

# --- "Data" sheet: append two new weekly observations ---
$wsData = $wb.Worksheets.Item("Data")

# Row 98 - copy formatting from the last existing data row (97), then set new values
$wsData.Range("A97:B97").Copy($wsData.Range("A98:B98"))
$wsData.Cells.Item(98, 1).Value = 45147
$wsData.Cells.Item(98, 2).Value = 444.594

# Row 99 - same treatment
$wsData.Range("A97:B97").Copy($wsData.Range("A99:B99"))
$wsData.Cells.Item(99, 1).Value = 45154
$wsData.Cells.Item(99, 2).Value = 435.795

# --- "SeriesInfo" sheet: refresh the metadata pulled from the FRED API ---
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

# Keep these as plain text (not auto-converted to dates) - mirror the original
# formatting by forcing Text entry, then restoring the default "Normal" style
# so the cell's look matches the unstyled source cell.
$wsInfo.Cells.Item(3, 2).NumberFormat = "@"
$wsInfo.Cells.Item(3, 2).Value = "2023-08-22"
$wsInfo.Cells.Item(3, 2).Style = "Normal"

$wsInfo.Cells.Item(4, 2).NumberFormat = "@"
$wsInfo.Cells.Item(4, 2).Value = "2023-08-22"
$wsInfo.Cells.Item(4, 2).Style = "Normal"

$wsInfo.Cells.Item(7, 2).NumberFormat = "@"
$wsInfo.Cells.Item(7, 2).Value = "2023-08-16"
$wsInfo.Cells.Item(7, 2).Style = "Normal"

$wsInfo.Cells.Item(14, 2).NumberFormat = "@"
$wsInfo.Cells.Item(14, 2).Value = "2023-08-17 15:35:18-05"
$wsInfo.Cells.Item(14, 2).Style = "Normal"
